# Edit script: trims the deck from 24 to 19 slides.
#  - Drops 4 slides whose content is superseded / no longer wanted:
#      "Contact and Repository", "Simple Architecture (Easy View)",
#      "Model Execution Mode", "Implementation Deliverables"
#  - Merges the two screenshot slides ("Frontend UI Screenshot" and
#    "Backend API Screenshot") into a single "Application Screenshots
#    (Frontend + Backend)" slide showing both images side by side,
#    and places it last in the deck.
#
# Deleting everything else simply shifts the remaining slides into
# their correct final positions/titles without any further text edits.

$p = $ppt.ActivePresentation

# --- Locate the two screenshot slides (by their current, pre-edit
#     positions) before any deletions happen, so we can copy their
#     shapes into the new merged slide. ---
$frontendSlide = $p.Slides.Item(18)
$backendSlide = $p.Slides.Item(19)

# --- Build the new, merged screenshot slide at the very end of the
#     deck, reusing the "Title Only" layout of the screenshot slides. ---
$layout = $frontendSlide.CustomLayout
$mergedSlide = $p.Slides.AddSlide($p.Slides.Count + 1, $layout)
$mergedSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Application Screenshots (Frontend + Backend)"

# Frontend screenshot (left half)
$frontendSlide.Shapes.Item(2).Copy()
$frontPic = $mergedSlide.Shapes.Paste()
$frontPic.Left = 36#
$frontPic.Top = 93.6
$frontPic.Width = 446.4
$frontPic.Height = 345.6

# Backend screenshot (right half)
$backendSlide.Shapes.Item(2).Copy()
$backPic = $mergedSlide.Shapes.Paste()
$backPic.Name = "Picture 3"
$backPic.Left = 475.2
$backPic.Top = 93.6
$backPic.Width = 446.4
$backPic.Height = 345.6

# Caption textbox - copied from the frontend slide's caption so it keeps
# the same styling (no-fill, autosize, 14pt default run size).
$frontendSlide.Shapes.Item(3).Copy()
$caption = $mergedSlide.Shapes.Paste()
$caption.Name = "TextBox 4"
$caption.TextFrame.TextRange.Text = "Left: Frontend UI demo | Right: Backend FastAPI Swagger docs"
$caption.Left = 43.2
$caption.Top = 453.6
$caption.Width = 871.2
$caption.Height = 50.4

# --- Now remove the slides that are no longer wanted. Delete from the
#     highest index down so earlier indices stay valid. ---
$backendSlide.Delete()
$frontendSlide.Delete()
$p.Slides.Item(17).Delete()   # "Implementation Deliverables"
$p.Slides.Item(14).Delete()   # "Model Execution Mode"
$p.Slides.Item(11).Delete()   # "Simple Architecture (Easy View)"
$p.Slides.Item(3).Delete()    # "Contact and Repository"

Write-Host "Final slide count:" $p.Slides.Count
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    Write-Host $i ":" $p.Slides.Item($i).Shapes.Item(1).TextFrame.TextRange.Text
}
